$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the R10 rule (row 8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 as the active cell, matching the sheet view selection in the saved file
$ws.Range("E8").Select()
